{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"GDP growth in 2019 decreased\", \"GDP growth in 2019 slowed down\"],\n  [\"GDP growth decreased from 2.5% in 2018 to 1.3% in 2019. On the demand side, private consumption, accounting for 50.4% of GDP, contributed the most with 1.3 percentage points (pp).\", \"GDP growth slowed down from 2.5% in 2018 to 1.3% in 2019. On the demand side, private consumption, accounting for 50.4% of GDP, contributed the most with 1.3 percentage points (pp).\"],\n  [\" Gross capital formation gave 0.8pp.\", \" Gross capital formation added 0.8pp.\"],\n  [\" Government expenditure shared 0.4pp.\", \" Government expenditure added 0.4pp.\"],\n  [\" On the other hand, net exports cut 1.4pp from growth.\", \" On the other hand, net exports shaved 1.4pp from growth.\"],\n  [\" Private consumption jumped by 2.5%.\", \" Private consumption increased by 2.5%.\"],\n  [\" Government expenditure picked up by 2.2%.\", \" Government expenditure jumped by 2.2%.\"],\n  [\"Services jumped by the largest edge at 1.6% annual growth.\", \"Services expanded by the largest edge at 1.6% annual growth.\"],\n  [\" Industry (including construction) grew by 1.0%.\", \" Industry (including construction) expanded by 1.0%.\"],\n  [\" Agriculture picked up by 0.6%.\", \" Agriculture increased by 0.6%.\"],\n  [\"Unemployment improved; inflation increased\", \"Unemployment improved; inflation jumped\"],\n  [\"Unemployment rate improved from 4.8% in 2018 to 4.5% in 2019. Consequently, inflation increased from 2.9% to 4.5%. At the end of the year, the central bank set the policy rate at 6.25%.\", \"Unemployment rate improved from 4.8% in 2018 to 4.5% in 2019. Consequently, inflation jumped from 2.9% to 4.5%. At the end of the year, the central bank set the policy rate at 6.25%.\"],\n  [\"Output plunged by 8.1% year-on-year in Q2 of 2020. Growth in overall economic activity worsened from 1.7% in the previous quarter. Net exports increased by the biggest margin at 227.5% annual growth.\", \"Output plunged by 8.1% year-on-year in Q2 of 2020. Growth in overall economic activity worsened from 1.7% in the previous quarter. Net exports expanded by the biggest margin at 227.5% annual growth.\"],\n  [\" Government expenditure picked up by 1.6%.\", \" Government expenditure increased by 1.6%.\"],\n  [\" On the other hand, private consumption and gross capital formation decreased by 22.2% and 6.2%, respectively.\", \" On the other hand, private consumption and gross capital formation contracted by 22.2% and 6.2%, respectively.\"],\n  [\"Retail sales contracted by 2.5% year-on-year in October of 2020. Growth in the retail sector increased from a contraction of 3.1% in September, reflecting increased trade activity.\", \"Retail sales contracted by 2.5% year-on-year in October of 2020. Growth in the retail sector jumped from a contraction of 3.1% in September, reflecting increased trade activity.\"],\n  [\"Meanwhile, consumer confidence index was in the negative territory at -22.0 points in Q3 of 2020. Confidence improved from -30.0 points in the previous quarter. Expectations of consumers about the general economic situation in the next 12 months turned pessimistic at -15.0 points from -20.0 points over the same period, reflecting improved consumer sentiments.\", \"Meanwhile, consumer confidence index was in the negative territory at -22.0 points in Q3 of 2020. Confidence increased from -30.0 points in the previous quarter. Expectations of consumers about the general economic situation in the next 12 months turned pessimistic at -15.0 points from -20.0 points over the same period, reflecting improved consumer sentiments.\"],\n  [\"Industrial production shrank by 5.6% year-on-year in October, a decrease from -5.1% growth in the previous month. Looking at the details, growth in manufacturing worsened to -3.5% from -0.7%, while mining and quarrying output growth improved to -8.8% from -9.5%. Meanwhile, growth in water supply, sewerage, waste management & remediation jumped to 1.9% from -0.9%, while electricity, gas, steam and air conditioning supply output growth declined to -3.2% from -2.6%\", \"Industrial production shrank by 5.6% year-on-year in October, a decrease from -5.1% growth in the previous month. Looking at the details, growth in manufacturing worsened to -3.5% from -0.7%, while mining and quarrying output growth increased to -8.8% from -9.5%. Meanwhile, growth in water supply, sewerage, waste management & remediation increased to 1.9% from -0.9%, while electricity, gas, steam and air conditioning supply output growth declined to -3.2% from -2.6%\"],\n  [\"Inflation rose\", \"Inflation worsened\"],\n  [\"Overall inflation rose to 4.0% year-on-year in October from 3.7% in the previous month. Prices for food products worsened to 4.8% from 4.3%, while housing, rent, water, electricity, gas & other fuels jumped to 3.4% from 3.3%. Transportation slowed down to -0.1% from 0.5%, while communication slowed down to 4.3% from 5.2%. Meanwhile, prices for health/medical care rose to 4.2% from 4.1%, recreation rose to 2.3% from 2.0%, and education improved to 1.9% from 2.1%\", \"Overall inflation worsened to 4.0% year-on-year in October from 3.7% in the previous month. Prices for food products worsened to 4.8% from 4.3%, while housing, rent, water, electricity, gas & other fuels jumped to 3.4% from 3.3%. Transportation slowed down to -0.1% from 0.5%, while communication slowed down to 4.3% from 5.2%. Meanwhile, prices for health/medical care worsened to 4.2% from 4.1%, recreation rose to 2.3% from 2.0%, and education slowed down to 1.9% from 2.1%\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"GDP growth in 2019 decreased\", \"GDP growth in 2019 slowed down\")\n  ,@(\"GDP growth decreased from 2.5% in 2018 to 1.3% in 2019. On the demand side, private consumption, accounting for 50.4% of GDP, contributed the most with 1.3 percentage points (pp).\", \"GDP growth slowed down from 2.5% in 2018 to 1.3% in 2019. On the demand side, private consumption, accounting for 50.4% of GDP, contributed the most with 1.3 percentage points (pp).\")\n  ,@(\" Gross capital formation gave 0.8pp.\", \" Gross capital formation added 0.8pp.\")\n  ,@(\" Government expenditure shared 0.4pp.\", \" Government expenditure added 0.4pp.\")\n  ,@(\" On the other hand, net exports cut 1.4pp from growth.\", \" On the other hand, net exports shaved 1.4pp from growth.\")\n  ,@(\" Private consumption jumped by 2.5%.\", \" Private consumption increased by 2.5%.\")\n  ,@(\" Government expenditure picked up by 2.2%.\", \" Government expenditure jumped by 2.2%.\")\n  ,@(\"Services jumped by the largest edge at 1.6% annual growth.\", \"Services expanded by the largest edge at 1.6% annual growth.\")\n  ,@(\" Industry (including construction) grew by 1.0%.\", \" Industry (including construction) expanded by 1.0%.\")\n  ,@(\" Agriculture picked up by 0.6%.\", \" Agriculture increased by 0.6%.\")\n  ,@(\"Unemployment improved; inflation increased\", \"Unemployment improved; inflation jumped\")\n  ,@(\"Unemployment rate improved from 4.8% in 2018 to 4.5% in 2019. Consequently, inflation increased from 2.9% to 4.5%. At the end of the year, the central bank set the policy rate at 6.25%.\", \"Unemployment rate improved from 4.8% in 2018 to 4.5% in 2019. Consequently, inflation jumped from 2.9% to 4.5%. At the end of the year, the central bank set the policy rate at 6.25%.\")\n  ,@(\"Output plunged by 8.1% year-on-year in Q2 of 2020. Growth in overall economic activity worsened from 1.7% in the previous quarter. Net exports increased by the biggest margin at 227.5% annual growth.\", \"Output plunged by 8.1% year-on-year in Q2 of 2020. Growth in overall economic activity worsened from 1.7% in the previous quarter. Net exports expanded by the biggest margin at 227.5% annual growth.\")\n  ,@(\" Government expenditure picked up by 1.6%.\", \" Government expenditure increased by 1.6%.\")\n  ,@(\" On the other hand, private consumption and gross capital formation decreased by 22.2% and 6.2%, respectively.\", \" On the other hand, private consumption and gross capital formation contracted by 22.2% and 6.2%, respectively.\")\n  ,@(\"Retail sales contracted by 2.5% year-on-year in October of 2020. Growth in the retail sector increased from a contraction of 3.1% in September, reflecting increased trade activity.\", \"Retail sales contracted by 2.5% year-on-year in October of 2020. Growth in the retail sector jumped from a contraction of 3.1% in September, reflecting increased trade activity.\")\n  ,@(\"Meanwhile, consumer confidence index was in the negative territory at -22.0 points in Q3 of 2020. Confidence improved from -30.0 points in the previous quarter. Expectations of consumers about the general economic situation in the next 12 months turned pessimistic at -15.0 points from -20.0 points over the same period, reflecting improved consumer sentiments.\", \"Meanwhile, consumer confidence index was in the negative territory at -22.0 points in Q3 of 2020. Confidence increased from -30.0 points in the previous quarter. Expectations of consumers about the general economic situation in the next 12 months turned pessimistic at -15.0 points from -20.0 points over the same period, reflecting improved consumer sentiments.\")\n  ,@(\"Industrial production shrank by 5.6% year-on-year in October, a decrease from -5.1% growth in the previous month. Looking at the details, growth in manufacturing worsened to -3.5% from -0.7%, while mining and quarrying output growth improved to -8.8% from -9.5%. Meanwhile, growth in water supply, sewerage, waste management & remediation jumped to 1.9% from -0.9%, while electricity, gas, steam and air conditioning supply output growth declined to -3.2% from -2.6%\", \"Industrial production shrank by 5.6% year-on-year in October, a decrease from -5.1% growth in the previous month. Looking at the details, growth in manufacturing worsened to -3.5% from -0.7%, while mining and quarrying output growth increased to -8.8% from -9.5%. Meanwhile, growth in water supply, sewerage, waste management & remediation increased to 1.9% from -0.9%, while electricity, gas, steam and air conditioning supply output growth declined to -3.2% from -2.6%\")\n  ,@(\"Inflation rose\", \"Inflation worsened\")\n  ,@(\"Overall inflation rose to 4.0% year-on-year in October from 3.7% in the previous month. Prices for food products worsened to 4.8% from 4.3%, while housing, rent, water, electricity, gas & other fuels jumped to 3.4% from 3.3%. Transportation slowed down to -0.1% from 0.5%, while communication slowed down to 4.3% from 5.2%. Meanwhile, prices for health/medical care rose to 4.2% from 4.1%, recreation rose to 2.3% from 2.0%, and education improved to 1.9% from 2.1%\", \"Overall inflation worsened to 4.0% year-on-year in October from 3.7% in the previous month. Prices for food products worsened to 4.8% from 4.3%, while housing, rent, water, electricity, gas & other fuels jumped to 3.4% from 3.3%. Transportation slowed down to -0.1% from 0.5%, while communication slowed down to 4.3% from 5.2%. Meanwhile, prices for health/medical care worsened to 4.2% from 4.1%, recreation rose to 2.3% from 2.0%, and education slowed down to 1.9% from 2.1%\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Text = $oldText\n  $rng.Find.Replacement.Text = $newText\n  $rng.Find.Forward = $true\n  $rng.Find.Wrap = 1\n  $rng.Find.MatchCase = $true\n  $rng.Find.MatchWholeWord = $false\n  $rng.Find.MatchWildcards = $false\n  $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n}"}
